$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 had "Bruno  da Silva" sent to Dóris' e-mail by mistake; correct the
# name to match the e-mail owner and drop the stale "ENVIADO" status so it
# gets re-sent via the corporate IMAP/SMTP.
$ws.Range("A4").Value = "Dóris Andressa Moura Luvizute"
$ws.Range("C4").ClearContents()
